$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format from the neighboring "sum" header (G1) onto the
# new H1 header cell so it keeps the same bold/border/alignment style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column value for the single data row.
$ws.Range("H2").Value = 0
